# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are stored as plain text, so a leading
# apostrophe is used for purely-numeric Price values to stop Excel's COM layer
# from auto-converting the assigned string into a numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '60.970.94'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.383.58'
$ws.Range("E3").Value = '  +0.02%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.00%  '

# Row 5: BNB
$ws.Range("D5").Value = '''571.12'
$ws.Range("E5").Value = '  -0.07%  '

# Row 6: Solana
$ws.Range("D6").Value = '''141.84'
$ws.Range("E6").Value = '  +0.21%  '

# Row 8: XRP
$ws.Range("E8").Value = '  +0.19%  '

# Row 9: Toncoin
$ws.Range("D9").Value = '''7.66'

# Row 10: Dogecoin
$ws.Range("E10").Value = '  -0.81%  '

# Row 11: Cardano
$ws.Range("E11").Value = '  -1.56%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '3.961.94'
$ws.Range("E12").Value = '  +0.04%  '

# Row 13: TRON
$ws.Range("E13").Value = '  +1.99%  '

# Row 14: Avalanche
$ws.Range("E14").Value = '  -0.87%  '

# Row 15: ShibaInu
$ws.Range("E15").Value = '  +0.20%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '3.386.19'
$ws.Range("E16").Value = '  -0.42%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '61.065.31'
$ws.Range("E17").Value = '  +0.29%  '

# Row 18: Polkadot
$ws.Range("D18").Value = '''6.09'
$ws.Range("E18").Value = '  -2.98%  '

# Row 19: Chainlink
$ws.Range("D19").Value = '''13.64'
$ws.Range("E19").Value = '  -3.69%  '

# Row 20: Uniswap
$ws.Range("E20").Value = '  -1.08%  '

# Row 21: BitcoinCash
$ws.Range("D21").Value = '''382.75'
$ws.Range("E21").Value = '  -1.52%  '

# Row 22: Litecoin
$ws.Range("D22").Value = '''75.25'

# Row 23: Polygon
$ws.Range("E23").Value = '  -1.47%  '

# Row 24: Dai
$ws.Range("E24").Value = '  +0.57%  '

# Row 25: PEPE
$ws.Range("E25").Value = '  -1.87%  '

# Row 26: WrappedeETH
$ws.Range("D26").Value = '3.522.42'
$ws.Range("E26").Value = '  -0.06%  '

# Row 27: Kaspa
$ws.Range("E27").Value = '  +2.19%  '

# Row 28: Binance-PegBSC-USD
$ws.Range("E28").Value = '  -0.06%  '

# Row 29: RenderToken
$ws.Range("D29").Value = '''7.23'
$ws.Range("E29").Value = '  -2.48%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = '''7.96'
$ws.Range("E30").Value = '  -1.49%  '

# Row 31: PancakeSwap
$ws.Range("D31").Value = '''2.14'
$ws.Range("E31").Value = '  -1.29%  '

# Row 33: Fetch.AI
$ws.Range("D33").Value = '''1.39'
$ws.Range("E33").Value = '  -5.07%  '

# Row 34: EthereumClassic
$ws.Range("D34").Value = '''23.18'
$ws.Range("E34").Value = '  -2.37%  '

# Row 35: Aptos
$ws.Range("E35").Value = '  +0.27%  '

# Row 36: Monero
$ws.Range("D36").Value = '''166.46'
$ws.Range("E36").Value = '  -0.41%  '

# Row 37: RenzoRestakedETH
$ws.Range("D37").Value = '3.415.94'
$ws.Range("E37").Value = '  +0.10%  '

# Row 38: NEARProtocol
$ws.Range("D38").Value = '''4.97'
$ws.Range("E38").Value = '  -1.25%  '

# Row 39: ImmutableX
$ws.Range("E39").Value = '  -2.74%  '

# Row 40: Hedera
$ws.Range("E40").Value = '  -1.70%  '

# Row 41: EnergySwap
$ws.Range("D41").Value = '''26.65'
$ws.Range("E41").Value = '  -1.23%  '

# Row 42: FirstDigitalUSD
$ws.Range("E42").Value = '  +0.01%  '

# Row 43: Mantle
$ws.Range("E43").Value = '  -0.59%  '

# Row 44: Filecoin
$ws.Range("D44").Value = '''4.37'
$ws.Range("E44").Value = '  -2.07%  '

# Row 45: Stacks
$ws.Range("E45").Value = '  -1.82%  '

# Row 46: ONDO
$ws.Range("E46").Value = '  -0.14%  '

# Row 47: Maker
$ws.Range("D47").Value = '2.445.39'
$ws.Range("E47").Value = '  -3.70%  '

# Row 48: InjectiveProtocol
$ws.Range("D48").Value = '''22.93'
$ws.Range("E48").Value = '  -0.22%  '

# Row 49: Cosmos
$ws.Range("E49").Value = '  -2.13%  '

# Row 50: dogwifhat
$ws.Range("E50").Value = '  +10.42%  '

# Row 51: VeChain
$ws.Range("E51").Value = '  -0.08%  '
